$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")

# Row 34 (key 33): Olympus Save Room
$ws.Cells.Item(34, 2).Value = "0x1f"
$ws.Cells.Item(34, 3).Value = "0x3"
$ws.Cells.Item(34, 4).Value = "OlympusSave"
$ws.Cells.Item(34, 5).Value = "Save Room"

# Row 35 (key 34): Olympus Layer 2
$ws.Cells.Item(35, 2).Value = "0x2"
$ws.Cells.Item(35, 3).Value = "0x3"
$ws.Cells.Item(35, 4).Value = "OlympusLayer2"
$ws.Cells.Item(35, 5).Value = "Layer 2"

# Row 36 (key 35): Olympus Layer 3
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "0x3"
$ws.Cells.Item(36, 3).Value = "0x3"
$ws.Cells.Item(36, 4).Value = "OlympusLayer3"
$ws.Cells.Item(36, 5).Value = "Layer 3"
$ws.Cells.Item(36, 6).Formula = '=_xlfn.CONCAT( ,A36,": { ""worldId"": ",C36,", ""name"": """,D36,""", ""display"": """,E36,""", ""areaId"": ",B36,", },")'

# Row 37 (key 36): Olympus Layer 4
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "0x4"
$ws.Cells.Item(37, 3).Value = "0x3"
$ws.Cells.Item(37, 4).Value = "OlympusLayer4"
$ws.Cells.Item(37, 5).Value = "Layer 4"
$ws.Cells.Item(37, 6).Formula = '=_xlfn.CONCAT( ,A37,": { ""worldId"": ",C37,", ""name"": """,D37,""", ""display"": """,E37,""", ""areaId"": ",B37,", },")'

# Row 38 (key 37): Olympus Layer 5
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "0x5"
$ws.Cells.Item(38, 3).Value = "0x3"
$ws.Cells.Item(38, 4).Value = "OlympusLayer5"
$ws.Cells.Item(38, 5).Value = "Layer 5"
$ws.Cells.Item(38, 6).Formula = '=_xlfn.CONCAT( ,A38,": { ""worldId"": ",C38,", ""name"": """,D38,""", ""display"": """,E38,""", ""areaId"": ",B38,", },")'

# Row 39 (key 38): Olympus Layer 6
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "0x6"
$ws.Cells.Item(39, 3).Value = "0x3"
$ws.Cells.Item(39, 4).Value = "OlympusLayer6"
$ws.Cells.Item(39, 5).Value = "Layer 6"
$ws.Cells.Item(39, 6).Formula = '=_xlfn.CONCAT( ,A39,": { ""worldId"": ",C39,", ""name"": """,D39,""", ""display"": """,E39,""", ""areaId"": ",B39,", },")'

# Rows 40-43 (keys 39-42): blank template rows, only key + formula
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 6).Formula = '=_xlfn.CONCAT( ,A40,": { ""worldId"": ",C40,", ""name"": """,D40,""", ""display"": """,E40,""", ""areaId"": ",B40,", },")'

$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 6).Formula = '=_xlfn.CONCAT( ,A41,": { ""worldId"": ",C41,", ""name"": """,D41,""", ""display"": """,E41,""", ""areaId"": ",B41,", },")'

$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 6).Formula = '=_xlfn.CONCAT( ,A42,": { ""worldId"": ",C42,", ""name"": """,D42,""", ""display"": """,E42,""", ""areaId"": ",B42,", },")'

$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 6).Formula = '=_xlfn.CONCAT( ,A43,": { ""worldId"": ",C43,", ""name"": """,D43,""", ""display"": """,E43,""", ""areaId"": ",B43,", },")'

# Match the author's final selection on the sheet
$ws.Activate()
$ws.Range("E41").Select()
